# Update "base_dados_pesquisa_PO" worksheet:
#  - add 4 new columns (F:I) with header labels in row 1, styled like the
#    existing header cells (A1:E1)
#  - append a new data row (row 9) with a survey response that includes
#    values for the new columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:I1, copying the header style from A1 ---
$ws.Range("A1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "Paineis Utilizados"
$ws.Range("G1").Value = "Painel com Feedback"
$ws.Range("H1").Value = "Comentário do Painel"
$ws.Range("I1").Value = "Data/Hora do Envio"

# --- New data row 9 ---
$ws.Range("A9").Value = "cleiton.souza@mrv.com.br"
$ws.Range("B9").Value = "Planilha geral - teste"
$ws.Range("C9").Value = "Painel Power BI"
$ws.Range("D9").Value = 4
$ws.Range("F9").Value = "Painel Análises Forecast de Produção - PLNESROBR009; Painel Operações - Planejamento e Controle - PLNESROBR010; Painel Produção Produtividade e MO - PLNESROBR005; Painel do Portifólio - Planejamento da Produção - PLNESROBR004; PAP - Dossiê"
$ws.Range("G9").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004"
$ws.Range("H9").Value = "Muito completo" + [char]10
$ws.Range("I9").Value = "2025-05-19 19:02:08"

Write-Host "Updated worksheet with new columns and row 9"
